$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The former column C (header "image2_url" plus hyperlinked frame-image URLs)
# needs to end up as column B once the old column B ("video_url" /
# "HitX_video_url_data") is removed. Capture the existing hyperlink
# targets (in their current collection order) before doing any structural
# edits, since this runtime does not auto-shift Hyperlink ranges when
# cells move.
$hyperlinkRows = @()
$hyperlinkAddrs = @()
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    $rowNum = [int]($addr -replace '[^0-9]', '')
    $hyperlinkRows += $rowNum
    $hyperlinkAddrs += $hl.Address
}

# Drop the stale hyperlink definitions; they will be recreated after the
# column shift so they point at the correct (new) cells.
$ws.Hyperlinks.Delete()

# Delete column B entirely; this shifts the old column C left into column B.
$ws.Range("B:B").Delete()

# Recreate the hyperlinks (same order as originally enumerated) anchored on
# the new column B cells.
for ($i = 0; $i -lt $hyperlinkRows.Count; $i++) {
    $target = $ws.Cells.Item($hyperlinkRows[$i], 2)
    $ws.Hyperlinks.Add($target, $hyperlinkAddrs[$i]) | Out-Null
}

# Adding a hyperlink re-stamps the cell with a fresh "Hyperlink" style entry;
# reapply the named style so the cells keep using the workbook's existing
# Hyperlink cell style (matching how the sheet looked pre-edit for column C).
$ws.Range("B2:B19").Style = "Hyperlink"

# Update the view selection to match the post-edit state.
$ws.Range("B1").Select()
